$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '66.555.95'
Set-TextValue 'E2' '  +4.42%  '
Set-TextValue 'D3' '3.493.68'
Set-TextValue 'E3' '  +2.24%  '
Set-TextValue 'E4' '  -0.03%  '
Set-TextValue 'D5' '590.53'
Set-TextValue 'E5' '  +3.74%  '
Set-TextValue 'D6' '170.08'
Set-TextValue 'E6' '  +8.35%  '
Set-TextValue 'D8' '3.495.89'
Set-TextValue 'E8' '  +2.14%  '
Set-TextValue 'D9' '0.584'
Set-TextValue 'E9' '  +2.31%  '
Set-TextValue 'E10' '  +0.26%  '
Set-TextValue 'D11' '0.125'
Set-TextValue 'E11' '  +4.31%  '
Set-TextValue 'E12' '  +2.22%  '
Set-TextValue 'D13' '4.099.16'
Set-TextValue 'E13' '  +2.25%  '
Set-TextValue 'E14' '  +0.57%  '
Set-TextValue 'E15' '  +3.30%  '
Set-TextValue 'D16' '66.501.84'
Set-TextValue 'E16' '  +4.18%  '
Set-TextValue 'E17' '  +1.85%  '
Set-TextValue 'D18' '3.494.59'
Set-TextValue 'E18' '  +2.62%  '
Set-TextValue 'D19' '6.29'
Set-TextValue 'E19' '  +2.90%  '
Set-TextValue 'E20' '  +3.07%  '
Set-TextValue 'D21' '388.39'
Set-TextValue 'E21' '  +1.60%  '
Set-TextValue 'E22' '  +2.73%  '
Set-TextValue 'E23' '  +2.50%  '
Set-TextValue 'D24' '1.00'
Set-TextValue 'E24' '  +0.07%  '
Set-TextValue 'E25' '  +1.73%  '
Set-TextValue 'D26' '0.0000124'
Set-TextValue 'E26' '  +7.37%  '
Set-TextValue 'D27' '10.08'
Set-TextValue 'E27' '  +4.19%  '
Set-TextValue 'E28' '  +1.76%  '
Set-TextValue 'D29' '0.999'
Set-TextValue 'E29' '  -0.12%  '
Set-TextValue 'D30' '6.38'
Set-TextValue 'E30' '  +4.69%  '
Set-TextValue 'E31' '  +6.21%  '
Set-TextValue 'E32' '  +3.54%  '
Set-TextValue 'D33' '23.49'
Set-TextValue 'E33' '  +2.23%  '
Set-TextValue 'D34' '7.39'
Set-TextValue 'E34' '  +5.57%  '
Set-TextValue 'E35' '  +0.10%  '
Set-TextValue 'E36' '  +1.38%  '
Set-TextValue 'D37' '160.95'
Set-TextValue 'E37' '  -0.15%  '
Set-TextValue 'D38' '0.902'
Set-TextValue 'E38' '  +8.07%  '
Set-TextValue 'E39' '  +5.70%  '
Set-TextValue 'D40' '0.0745'
Set-TextValue 'E40' '  +2.92%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D41' '6.72'
Set-TextValue 'E41' '  +5.24%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D42' '26.49'
Set-TextValue 'E42' '  +1.04%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D43' '4.60'
Set-TextValue 'E43' '  +4.29%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D44' '27.04'
Set-TextValue 'E44' '  +5.01%  '
Set-TextValue 'D45' '2.812.27'
Set-TextValue 'E45' '  +0.09%  '
Set-TextValue 'D46' '43.50'
Set-TextValue 'E46' '  +1.60%  '
$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D47' '2.56'
Set-TextValue 'E47' '  +10.48%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D48' '0.0315'
Set-TextValue 'E48' '  +3.01%  '
$ws.Range('B49').Value = 'Bittensor'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D49' '355.89'
Set-TextValue 'E49' '  +8.61%  '
Set-TextValue 'E50' '  +5.86%  '
Set-TextValue 'D51' '32.89'
Set-TextValue 'E51' '  +9.59%  '
